$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testData")

# Replace the 10 existing customer id rows (A13:A22) with 8 new customer ids (A13:A20)
$ws.Range("A13").Value = "cus_Gcpad6iDgpeoQg"
$ws.Range("A14").Value = "cus_Gcpa6eGWaWaheA"
$ws.Range("A15").Value = "cus_GcpT6DKbq2kFUu"
$ws.Range("A16").Value = "cus_GcpTnjCq9GTofk"
$ws.Range("A17").Value = "cus_GckBkHUVz22jOB"
$ws.Range("A18").Value = "cus_GckBYUUD1XZNaI"
$ws.Range("A19").Value = "cus_GckBH0AAVclMZA"
$ws.Range("A20").Value = "cus_GckBRoFgreGY6l"

# Remove the now-unused trailing rows 21 and 22 (delete bottom row first)
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(21).Delete()

# Re-select A13:A20 so the worksheet's saved selection matches the new data extent
$ws.Activate() | Out-Null
$ws.Range("A13:A20").Select() | Out-Null
